$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.95628266666667
$ws.Range("H2").Value = 86.868848
$ws.Range("I2").Value = 0.5491054194301004
$ws.Range("J2").Value = 0.5491054194301005
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 23.34377233333333
$ws.Range("N2").Value = 70.031317
$ws.Range("O2").Value = 0.6243713731385436
$ws.Range("P2").Value = 0.6243713731385436
$ws.Range("Q2").Value = 675.9488701903128
$ws.Range("R2").Value = 6083.539831712817
$ws.Range("S2").Value = 0.3428457047273877
$ws.Range("T2").Value = 0.3428457047273878

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.95628266666667
$ws.Range("H3").Value = 86.868848
$ws.Range("I3").Value = 0.5491054194301004
$ws.Range("J3").Value = 0.5491054194301005
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.314670666666666
$ws.Range("N3").Value = 24.944012
$ws.Range("O3").Value = 0.2223908915496236
$ws.Range("P3").Value = 0.2223908915496236
$ws.Range("Q3").Value = 240.7619541042418
$ws.Range("R3").Value = 2166.857586938176
$ws.Range("S3").Value = 0.12211604378179
$ws.Range("T3").Value = 0.12211604378179

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 28.95628266666667
$ws.Range("H4").Value = 86.868848
$ws.Range("I4").Value = 0.5491054194301004
$ws.Range("J4").Value = 0.5491054194301005
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02865466666666667
$ws.Range("N4").Value = 0.085964
$ws.Range("O4").Value = 0.0007664208388438813
$ws.Range("P4").Value = 0.0007664208388438813
$ws.Range("Q4").Value = 0.829732627719111
$ws.Range("R4").Value = 7.467593649472
$ws.Range("S4").Value = 0.0004208458361733388
$ws.Range("T4").Value = 0.0004208458361733389

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 28.95628266666667
$ws.Range("H5").Value = 86.868848
$ws.Range("I5").Value = 0.5491054194301004
$ws.Range("J5").Value = 0.5491054194301005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.700542666666667
$ws.Range("N5").Value = 17.101628
$ws.Range("O5").Value = 0.1524713144729888
$ws.Range("P5").Value = 0.1524713144729888
$ws.Range("Q5").Value = 165.0665248093938
$ws.Range("R5").Value = 1485.598723284544
$ws.Range("S5").Value = 0.08372282508474925
$ws.Range("T5").Value = 0.08372282508474928

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.691493
$ws.Range("H6").Value = 38.074479
$ws.Range("I6").Value = 0.2406720388519202
$ws.Range("J6").Value = 0.2406720388519202
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 23.34377233333333
$ws.Range("N6").Value = 70.031317
$ws.Range("O6").Value = 0.6243713731385436
$ws.Range("P6").Value = 0.6243713731385436
$ws.Range("Q6").Value = 296.2673231620936
$ws.Range("R6").Value = 2666.405908458843
$ws.Range("S6").Value = 0.1502687313740263
$ws.Range("T6").Value = 0.1502687313740263

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.691493
$ws.Range("H7").Value = 38.074479
$ws.Range("I7").Value = 0.2406720388519202
$ws.Range("J7").Value = 0.2406720388519202
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.314670666666666
$ws.Range("N7").Value = 24.944012
$ws.Range("O7").Value = 0.2223908915496236
$ws.Range("P7").Value = 0.2223908915496236
$ws.Range("Q7").Value = 105.5255845633053
$ws.Range("R7").Value = 949.7302610697479
$ws.Range("S7").Value = 0.05352326929134417
$ws.Range("T7").Value = 0.05352326929134417

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.691493
$ws.Range("H8").Value = 38.074479
$ws.Range("I8").Value = 0.2406720388519202
$ws.Range("J8").Value = 0.2406720388519202
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02865466666666667
$ws.Range("N8").Value = 0.085964
$ws.Range("O8").Value = 0.0007664208388438813
$ws.Range("P8").Value = 0.0007664208388438813
$ws.Range("Q8").Value = 0.3636705014173333
$ws.Range("R8").Value = 3.273034512756
$ws.Range("S8").Value = 0.0001844560659031558
$ws.Range("T8").Value = 0.0001844560659031558

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.691493
$ws.Range("H9").Value = 38.074479
$ws.Range("I9").Value = 0.2406720388519202
$ws.Range("J9").Value = 0.2406720388519202
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.700542666666667
$ws.Range("N9").Value = 17.101628
$ws.Range("O9").Value = 0.1524713144729888
$ws.Range("P9").Value = 0.1524713144729888
$ws.Range("Q9").Value = 72.34839735020134
$ws.Range("R9").Value = 651.135576151812
$ws.Range("S9").Value = 0.0366955821206465
$ws.Range("T9").Value = 0.0366955821206465

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Ephb4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4888703333333334
$ws.Range("H10").Value = 1.466611
$ws.Range("I10").Value = 0.009270573592685367
$ws.Range("J10").Value = 0.009270573592685367
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 23.34377233333333
$ws.Range("N10").Value = 70.031317
$ws.Range("O10").Value = 0.6243713731385436
$ws.Range("P10").Value = 0.6243713731385436
$ws.Range("Q10").Value = 11.41207776185411
$ws.Range("R10").Value = 102.708699856687
$ws.Range("S10").Value = 0.005788280763846884
$ws.Range("T10").Value = 0.005788280763846884

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Efnb2"
$ws.Range("C11").Value = "Ephb4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4888703333333334
$ws.Range("H11").Value = 1.466611
$ws.Range("I11").Value = 0.009270573592685367
$ws.Range("J11").Value = 0.009270573592685367
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.314670666666666
$ws.Range("N11").Value = 24.944012
$ws.Range("O11").Value = 0.2223908915496236
$ws.Range("P11").Value = 0.2223908915496236
$ws.Range("Q11").Value = 4.064795820370223
$ws.Range("R11").Value = 36.583162383332
$ws.Range("S11").Value = 0.002061691126453696
$ws.Range("T11").Value = 0.002061691126453696

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Efnb2"
$ws.Range("C12").Value = "Ephb4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4888703333333334
$ws.Range("H12").Value = 1.466611
$ws.Range("I12").Value = 0.009270573592685367
$ws.Range("J12").Value = 0.009270573592685367
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02865466666666667
$ws.Range("N12").Value = 0.085964
$ws.Range("O12").Value = 0.0007664208388438813
$ws.Range("P12").Value = 0.0007664208388438813
$ws.Range("Q12").Value = 0.01400841644488889
$ws.Range("R12").Value = 0.126075748004
$ws.Range("S12").Value = [double]"7.105160789469853E-06"
$ws.Range("T12").Value = [double]"7.105160789469853E-06"

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Efnb2"
$ws.Range("C13").Value = "Ephb4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4888703333333334
$ws.Range("H13").Value = 1.466611
$ws.Range("I13").Value = 0.009270573592685367
$ws.Range("J13").Value = 0.009270573592685367
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.700542666666667
$ws.Range("N13").Value = 17.101628
$ws.Range("O13").Value = 0.1524713144729888
$ws.Range("P13").Value = 0.1524713144729888
$ws.Range("Q13").Value = 2.786826193634222
$ws.Range("R13").Value = 25.081435742708
$ws.Range("S13").Value = 0.001413496541595316
$ws.Range("T13").Value = 0.001413496541595316

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Efnb2"
$ws.Range("C14").Value = "Ephb4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 10.59691233333333
$ws.Range("H14").Value = 31.790737
$ws.Range("I14").Value = 0.200951968125294
$ws.Range("J14").Value = 0.200951968125294
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 23.34377233333333
$ws.Range("N14").Value = 70.031317
$ws.Range("O14").Value = 0.6243713731385436
$ws.Range("P14").Value = 0.6243713731385436
$ws.Range("Q14").Value = 247.3719089456254
$ws.Range("R14").Value = 2226.347180510629
$ws.Range("S14").Value = 0.1254686562732827
$ws.Range("T14").Value = 0.1254686562732827

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Efnb2"
$ws.Range("C15").Value = "Ephb4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 10.59691233333333
$ws.Range("H15").Value = 31.790737
$ws.Range("I15").Value = 0.200951968125294
$ws.Range("J15").Value = 0.200951968125294
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.314670666666666
$ws.Range("N15").Value = 24.944012
$ws.Range("O15").Value = 0.2223908915496236
$ws.Range("P15").Value = 0.2223908915496236
$ws.Range("Q15").Value = 88.10983613520489
$ws.Range("R15").Value = 792.988525216844
$ws.Range("S15").Value = 0.04468988735003568
$ws.Range("T15").Value = 0.04468988735003568

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Efnb2"
$ws.Range("C16").Value = "Ephb4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 10.59691233333333
$ws.Range("H16").Value = 31.790737
$ws.Range("I16").Value = 0.200951968125294
$ws.Range("J16").Value = 0.200951968125294
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02865466666666667
$ws.Range("N16").Value = 0.085964
$ws.Range("O16").Value = 0.0007664208388438813
$ws.Range("P16").Value = 0.0007664208388438813
$ws.Range("Q16").Value = 0.3036509906075556
$ws.Range("R16").Value = 2.732858915468
$ws.Range("S16").Value = 0.0001540137759779167
$ws.Range("T16").Value = 0.0001540137759779167

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Efnb2"
$ws.Range("C17").Value = "Ephb4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 10.59691233333333
$ws.Range("H17").Value = 31.790737
$ws.Range("I17").Value = 0.200951968125294
$ws.Range("J17").Value = 0.200951968125294
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 5.700542666666667
$ws.Range("N17").Value = 17.101628
$ws.Range("O17").Value = 0.1524713144729888
$ws.Range("P17").Value = 0.1524713144729888
$ws.Range("Q17").Value = 60.40815089109289
$ws.Range("R17").Value = 543.6733580198361
$ws.Range("S17").Value = 0.03063941072599773
$ws.Range("T17").Value = 0.03063941072599773
